$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.996.16"
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").Value = "'2.238.75"
$ws.Range("E3").Value = '  -2.51%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("E5").Value = '  -3.82%  '
$ws.Range("D6").Value = "'86.49"
$ws.Range("E6").Value = '  +6.62%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = "'0.471"
$ws.Range("E9").Value = '  +0.59%  '
$ws.Range("E10").Value = '  +2.00%  '
$ws.Range("E11").Value = '  +7.43%  '
$ws.Range("D12").Value = "'47.12"
$ws.Range("E12").Value = '  -9.56%  '
$ws.Range("D13").Value = "'0.108"
$ws.Range("E13").Value = '  -0.61%  '
$ws.Range("D14").Value = "'6.41"
$ws.Range("E14").Value = '  +4.08%  '
$ws.Range("D15").Value = "'2.580.17"
$ws.Range("E15").Value = '  -2.63%  '
$ws.Range("D16").Value = "'14.21"
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = "'2.245.62"
$ws.Range("E17").Value = '  -2.54%  '
$ws.Range("D18").Value = "'0.732"
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").Value = "'39.921.47"
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("D20").Value = "'0.0₃0895"
$ws.Range("E20").Value = '  +2.31%  '
$ws.Range("D21").Value = "'5.81"
$ws.Range("E21").Value = '  -1.23%  '
$ws.Range("D22").Value = "'10.60"
$ws.Range("E22").Value = '  +4.96%  '
$ws.Range("D23").Value = "'65.59"
$ws.Range("E23").Value = '  -1.87%  '
$ws.Range("D24").Value = "'236.41"
$ws.Range("E24").Value = '  +2.96%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("E26").Value = '  -0.78%  '
$ws.Range("D27").Value = "'1.84"
$ws.Range("E27").Value = '  +4.47%  '
$ws.Range("D28").Value = "'23.04"
$ws.Range("E28").Value = '  +1.78%  '
$ws.Range("E29").Value = '  +1.55%  '
$ws.Range("E30").Value = '  +3.51%  '
$ws.Range("D31").Value = "'34.26"
$ws.Range("E31").Value = '  +6.69%  '
$ws.Range("D32").Value = "'154.36"
$ws.Range("E32").Value = '  +3.20%  '
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("D34").Value = "'4.88"
$ws.Range("E34").Value = '  -0.52%  '
$ws.Range("D35").Value = "'0.0712"
$ws.Range("E35").Value = '  +2.12%  '
$ws.Range("E36").Value = '  -2.10%  '
$ws.Range("D37").Value = "'16.60"
$ws.Range("E37").Value = '  +10.44%  '
$ws.Range("E38").Value = '  +0.46%  '
$ws.Range("D39").Value = "'0.100"
$ws.Range("D40").Value = "'2.71"
$ws.Range("E40").Value = '  +1.25%  '
$ws.Range("E41").Value = '  +2.75%  '
$ws.Range("D42").Value = "'3.80"
$ws.Range("E42").Value = '  +3.76%  '
$ws.Range("D43").Value = "'1.967.17"
$ws.Range("E44").Value = '  -2.30%  '
$ws.Range("E45").Value = '  +6.41%  '
$ws.Range("D46").Value = "'9.61"
$ws.Range("E46").Value = '  +5.65%  '
$ws.Range("D47").Value = "'16.37"
$ws.Range("E47").Value = '  -2.28%  '
$ws.Range("D48").Value = "'2.60"
$ws.Range("E48").Value = '  +1.00%  '
$ws.Range("D49").Value = "'2.451.98"
$ws.Range("E49").Value = '  -2.54%  '
$ws.Range("D50").Value = "'71.05"
$ws.Range("E50").Value = '  +4.85%  '
$ws.Range("D51").Value = "'1.47"
$ws.Range("E51").Value = '  +12.11%  '
